$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Append the "missing insurance" addendum after "baby_mom_at_birth.csv"
#    in the first bullet list (adds 5 new runs, including a Wingdings arrow
#    symbol character, to the existing paragraph).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(2)
$fullRange = $p1.Range

$pPr  = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>'
$run1 = '<w:r w:rsidRPr="002D488E"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>baby_mom_at_birth.csv</w:t></w:r>'
$run2 = '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
$run3 = '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>'
$run4 = '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> missing insurance (</w:t></w:r>'
$run5 = '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>baby_mom_at_birth_with_payer</w:t></w:r>'
$run6 = '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>.csv)</w:t></w:r>'

$para = '<w:p w14:paraId="6E01F970" w14:textId="5023BC26" w:rsidR="000D6EFF" w:rsidRPr="002D488E" w:rsidRDefault="000D6EFF" w:rsidP="002D488E">' + $pPr + $run1 + $run2 + $run3 + $run4 + $run5 + $run6 + '</w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $para + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$fullRange.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) Collapse filename runs that were previously split into two <w:r>
#    elements (name + ".csv") back into a single run/text node.
# ---------------------------------------------------------------------------
function Merge-Text($findText) {
    $range = $d.Content
    [void]$range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $findText, 2)
}

Merge-Text "mom_notes_at_delivery.csv"
Merge-Text "postnatal_baby_metadata.csv"
Merge-Text "postnatal_mom_metadata.csv"
Merge-Text "mom_notes_prenatal_visit.csv"
Merge-Text "subjects_clinical_notes_details_1-8.csv (infant) 8 files in total"
Merge-Text "subjects_clinical_notes_details_1-11.csv (mom) 11 files in total"

# ---------------------------------------------------------------------------
# 3) Shrink the page margins from 1 inch (1440 twips / 72 pt) down to
#    0.5 inch (720 twips / 36 pt) on all four sides.
# ---------------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $ps = $sec.PageSetup
    $ps.TopMargin = 36
    $ps.BottomMargin = 36
    $ps.LeftMargin = 36
    $ps.RightMargin = 36
}
